{"js": "// Swap \"Markdown,\" and \"Docusaurus,\" within the \"Documentation:\" skills\n// line so it reads \"...JSDoc, Docusaurus, Markdown, LaTeX\" instead of\n// \"...JSDoc, Markdown, Docusaurus, LaTeX\". Each word is its own run\n// (separated by a run that holds just the space), so we only rewrite the\n// text of the two word runs and leave the separating space run untouched\n// \u2014 this mirrors the target XML diff exactly (same run/rPr structure,\n// only the two <w:t> contents trade places).\n\nconst body = context.document.body;\n\nconst markdownResults = body.search(\"Markdown,\", { matchCase: true, matchWholeWord: false });\nmarkdownResults.load(\"text\");\nconst docusaurusResults = body.search(\"Docusaurus,\", { matchCase: true, matchWholeWord: false });\ndocusaurusResults.load(\"text\");\n\nawait context.sync();\n\nif (markdownResults.items.length === 0 || docusaurusResults.items.length === 0) {\n  throw new Error(\"Could not locate both 'Markdown,' and 'Docusaurus,' runs.\");\n}\n\nconst markdownRange = markdownResults.items[0];\nconst docusaurusRange = docusaurusResults.items[0];\n\n// Replace in place so each range keeps its own run formatting (rPr).\nmarkdownRange.insertText(\"Docusaurus,\", Word.InsertLocation.replace);\ndocusaurusRange.insertText(\"Markdown,\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Swap \"Markdown,\" and \"Docusaurus,\" within the \"Documentation:\" skills\n# line so it reads \"...JSDoc, Docusaurus, Markdown, LaTeX\" instead of\n# \"...JSDoc, Markdown, Docusaurus, LaTeX\". Each word lives in its own run\n# (with a separate run holding just the space between them), so we only\n# rewrite the text of the two word runs and leave the separating space\n# run untouched -- this mirrors the target XML diff exactly (same\n# run/rPr structure, only the two <w:t> contents trade places).\n#\n# A unique placeholder is used so the two words can be swapped without\n# one replacement accidentally colliding with the other.\n\n$d = $word.ActiveDocument\n$placeholder = \"DOCX_SWAP_PLACEHOLDER\"\n\n# 1) Find \"Markdown,\" and stash it behind a unique placeholder.\n$rngMarkdown = $d.Content\n$rngMarkdown.Find.ClearFormatting()\n$rngMarkdown.Find.Text = \"Markdown,\"\n$rngMarkdown.Find.MatchCase = $true\n$rngMarkdown.Find.MatchWholeWord = $false\n$rngMarkdown.Find.Execute() | Out-Null\n$rngMarkdown.Text = $placeholder\n\n# 2) Find \"Docusaurus,\" and rename it to \"Markdown,\".\n$rngDocusaurus = $d.Content\n$rngDocusaurus.Find.ClearFormatting()\n$rngDocusaurus.Find.Text = \"Docusaurus,\"\n$rngDocusaurus.Find.MatchCase = $true\n$rngDocusaurus.Find.MatchWholeWord = $false\n$rngDocusaurus.Find.Execute() | Out-Null\n$rngDocusaurus.Text = \"Markdown,\"\n\n# 3) Find the placeholder and rename it to \"Docusaurus,\".\n$rngPlaceholder = $d.Content\n$rngPlaceholder.Find.ClearFormatting()\n$rngPlaceholder.Find.Text = $placeholder\n$rngPlaceholder.Find.MatchCase = $true\n$rngPlaceholder.Find.MatchWholeWord = $false\n$rngPlaceholder.Find.Execute() | Out-Null\n$rngPlaceholder.Text = \"Docusaurus,\"\n"}
